# =====================================================================
# SEED 999 birth-results workbook: add BOUNDARY attack block (AK:AR)
# "add res boundary attack for SEED 999"
# =====================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Formatting for the new AK:AR header block (rows 1-2): copy from the
#    existing FGSM block (AC:AJ), an identical 8-column epsilon sweep
#    (0.01..0.20) with the bordered/bold/centered header style.
# ---------------------------------------------------------------------
$src = $ws.Range("AC1:AJ2")
$dst = $ws.Range("AK1:AR2")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Header text/labels for the new block.
#    "0.0x" values look numeric, so they are entered with a leading
#    apostrophe to force text storage (matching the other epsilon-sweep
#    header rows, which are text, not numbers), then formats are
#    re-applied on top so the quote-prefix style doesn't linger on the
#    cells themselves.
# ---------------------------------------------------------------------
$ws.Range("AK1").Value = "BOUNDARY"
$ws.Range("AK2").Value = "'0.01"
$ws.Range("AL2").Value = "'0.02"
$ws.Range("AM2").Value = "'0.03"
$ws.Range("AN2").Value = "'0.04"
$ws.Range("AO2").Value = "'0.05"
$ws.Range("AP2").Value = "'0.07"
$ws.Range("AQ2").Value = "'0.10"
$ws.Range("AR2").Value = "'0.20"

$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Merge the new header cell AK1:AR1 like the other attack headers.
# ---------------------------------------------------------------------
$ws.Range("AK1:AR1").Merge()

# ---------------------------------------------------------------------
# 4) Numeric data for the new BOUNDARY columns (rows 4,5,7,8,10,11 raw
#    metrics; rows 6,9,12 derived ratios) - plain numbers, default style.
# ---------------------------------------------------------------------
$ws.Range("AK4").Value = 431.0746963310241
$ws.Range("AL4").Value = 434.0184351921081
$ws.Range("AM4").Value = 442.0884850056966
$ws.Range("AN4").Value = 455.2186360168457
$ws.Range("AO4").Value = 445.0442343393962
$ws.Range("AP4").Value = 456.4610941314697
$ws.Range("AQ4").Value = 535.1031312561036
$ws.Range("AR4").Value = 675.1907266489665
$ws.Range("AK5").Value = 532.2268584325789
$ws.Range("AL5").Value = 540.3863566357007
$ws.Range("AM5").Value = 545.2303038161166
$ws.Range("AN5").Value = 562.6111922357381
$ws.Range("AO5").Value = 549.4706235084179
$ws.Range("AP5").Value = 553.2191127153668
$ws.Range("AQ5").Value = 652.1437605636816
$ws.Range("AR5").Value = 830.2693789194271
$ws.Range("AK6").Value = 0.9990560563949743
$ws.Range("AL6").Value = 0.9990173393555404
$ws.Range("AM6").Value = 0.9990006472173181
$ws.Range("AN6").Value = 0.9989340528622581
$ws.Range("AO6").Value = 0.9989856657458011
$ws.Range("AP6").Value = 0.9989978367756847
$ws.Range("AQ6").Value = 0.9986149362454411
$ws.Range("AR6").Value = 0.9975031909371861
$ws.Range("AK7").Value = 457.2896216392517
$ws.Range("AL7").Value = 460.9160689926147
$ws.Range("AM7").Value = 469.0771188735962
$ws.Range("AN7").Value = 478.7062649726868
$ws.Range("AO7").Value = 491.6180555788676
$ws.Range("AP7").Value = 491.0989400800069
$ws.Range("AQ7").Value = 522.3191157341004
$ws.Range("AR7").Value = 659.2025764338175
$ws.Range("AK8").Value = 570.3243797119386
$ws.Range("AL8").Value = 574.575567590014
$ws.Range("AM8").Value = 586.0018340183494
$ws.Range("AN8").Value = 594.7597998142642
$ws.Range("AO8").Value = 607.2216428211774
$ws.Range("AP8").Value = 614.203949683111
$ws.Range("AQ8").Value = 635.342043246599
$ws.Range("AR8").Value = 837.5006920164858
$ws.Range("AK9").Value = 0.9992637710186979
$ws.Range("AL9").Value = 0.9992512850217402
$ws.Range("AM9").Value = 0.9992258806774564
$ws.Range("AN9").Value = 0.9991819670309291
$ws.Range("AO9").Value = 0.999151201069204
$ws.Range("AP9").Value = 0.9990875000029596
$ws.Range("AQ9").Value = 0.9990212248333846
$ws.Range("AR9").Value = 0.9979624185505541
$ws.Range("AK10").Value = 273.088470808665
$ws.Range("AL10").Value = 276.7389280637105
$ws.Range("AM10").Value = 285.3536395581563
$ws.Range("AN10").Value = 302.9684701220194
$ws.Range("AO10").Value = 294.3688612937927
$ws.Range("AP10").Value = 364.7713779894511
$ws.Range("AQ10").Value = 385.0329274940491
$ws.Range("AR10").Value = 710.1054957135518
$ws.Range("AK11").Value = 402.3798521038682
$ws.Range("AL11").Value = 401.6029136052479
$ws.Range("AM11").Value = 406.9162247784251
$ws.Range("AN11").Value = 430.2524921809209
$ws.Range("AO11").Value = 419.6305782516142
$ws.Range("AP11").Value = 511.8494156363107
$ws.Range("AQ11").Value = 524.973840867762
$ws.Range("AR11").Value = 909.0024500478823
$ws.Range("AK12").Value = 0.9993614883184079
$ws.Range("AL12").Value = 0.999363617255234
$ws.Range("AM12").Value = 0.999345498506686
$ws.Range("AN12").Value = 0.999267625848628
$ws.Range("AO12").Value = 0.9993020313413635
$ws.Range("AP12").Value = 0.9989847573538428
$ws.Range("AQ12").Value = 0.9989066878463636
$ws.Range("AR12").Value = 0.9967095090081728

# ---------------------------------------------------------------------
# 5) Tiny ULP-level recomputation drift on pre-existing ratio rows
#    (6, 9, 12) caused by the upstream metric recompute after adding
#    the new attack column.
# ---------------------------------------------------------------------
$ws.Range("C6").Value = 0.9990572076996539
$ws.Range("F6").Value = 0.9990278856123233
$ws.Range("K6").Value = 0.9984632812530292
$ws.Range("M6").Value = 0.9990478684215719
$ws.Range("Q6").Value = 0.9989355060477946
$ws.Range("AA6").Value = 0.9916633721710616
$ws.Range("AC6").Value = 0.9987047030267312
$ws.Range("AD6").Value = 0.9982252584098521
$ws.Range("AF6").Value = 0.9969923053063227
$ws.Range("AH6").Value = 0.9932304717754413
$ws.Range("AI6").Value = 0.9892100130001565
$ws.Range("E9").Value = 0.9992695035361205
$ws.Range("G9").Value = 0.9992605939404476
$ws.Range("I9").Value = 0.9992278320109026
$ws.Range("M9").Value = 0.9992695035361205
$ws.Range("O9").Value = 0.9992627232274649
$ws.Range("Q9").Value = 0.9992280994810504
$ws.Range("V9").Value = 0.9989434509272842
$ws.Range("X9").Value = 0.9984823254235285
$ws.Range("Y9").Value = 0.9981943015910842
$ws.Range("AB9").Value = 0.9885143133645865
$ws.Range("AE9").Value = 0.9987383743891032
$ws.Range("AF9").Value = 0.9984910585475546
$ws.Range("D12").Value = 0.999254587623425
$ws.Range("H12").Value = 0.9992552063786594
$ws.Range("I12").Value = 0.999219932867908
$ws.Range("M12").Value = 0.9993890749560386
$ws.Range("N12").Value = 0.9993922548725127
$ws.Range("O12").Value = 0.9993837796828535
$ws.Range("P12").Value = 0.9993647643320814
$ws.Range("R12").Value = 0.9992539122027454
$ws.Range("T12").Value = 0.9983951346549725
$ws.Range("U12").Value = 0.9992081459484088
$ws.Range("V12").Value = 0.9989926318062669
$ws.Range("W12").Value = 0.9987276101008802
$ws.Range("X12").Value = 0.9984113088916848
$ws.Range("AB12").Value = 0.9840550853465725
$ws.Range("AE12").Value = 0.9985925687672162
$ws.Range("AF12").Value = 0.9982141732941139
$ws.Range("AG12").Value = 0.9977746498301803
